$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.021.90'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.978.08'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.22'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.23'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.426'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.10'
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.362'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.496.08'
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.91'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000160'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.019.43'
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.990.29'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.99'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.86'
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.95'
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.38'
$ws.Range("E21").Value = '  +2.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.491'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.47'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.115.57'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0891'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.41'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.97'
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.77'
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.11'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.82'
$ws.Range("E34").Value = '  -2.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.51'
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.71'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.21'
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.24'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0655'
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.012.76'
$ws.Range("E40").Value = '  +0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.50'
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.651'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.172.08'
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.34'
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.87'
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.920'
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.59'
$ws.Range("E50").Value = '  +2.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0849'
$ws.Range("E51").Value = '  -3.09%  '
